# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted for "Membrillo" (Vega Modelo de
# Temuco) at sheet row 82. All existing rows from 82 downward shift down by
# one (old row 82 -> new row 83, ..., old row 181 -> new row 182), and the
# freshly inserted row 82 is populated with this week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 82, pushing rows 82:181 down to 83:182
# (this also grows the sheet's used range / dimension from T181 to T182).
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly data point.
$ws.Range("A82").Value = 10
$ws.Range("B82").Value = 'Vega Modelo de Temuco'
$ws.Range("C82").Value = 'La Araucanía'
$ws.Range("D82").Value = 44740
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = 'Fruta'
$ws.Range("G82").Value = 100104
$ws.Range("H82").Value = 'Frutos de pepita'
$ws.Range("I82").Value = 100104003
$ws.Range("J82").Value = 'Membrillo'
$ws.Range("K82").Value = 'Champion'
$ws.Range("L82").Value = 'Primera'
$ws.Range("M82").Value = 100
$ws.Range("N82").Value = 10000
$ws.Range("O82").Value = 10000
$ws.Range("P82").Value = 10000
$ws.Range("Q82").Value = '$/bandeja 18 kilos granel'
$ws.Range("R82").Value = "Región de O'Higgins"
$ws.Range("S82").Value = 556
$ws.Range("T82").Value = 18
